$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, volume %, and price values that are
# not at risk of being auto-parsed as numbers).
$safeValues = [ordered]@{
    'D2' = '28.853.03'
    'E2' = '  +7.84%  '
    'D3' = '1.810.00'
    'E3' = '  +4.95%  '
    'E4' = '  +0.18%  '
    'E5' = '  +3.30%  '
    'E6' = '  +0.17%  '
    'E7' = '  +2.14%  '
    'E8' = '  +7.69%  '
    'E9' = '  +3.70%  '
    'D10' = '1.812.36'
    'E10' = '  +4.85%  '
    'E11' = '  +5.38%  '
    'E12' = '  +2.35%  '
    'E13' = '  +6.80%  '
    'E14' = '  +9.29%  '
    'E15' = '  +5.09%  '
    'D16' = '28.821.94'
    'E16' = '  +8.47%  '
    'E17' = '  +0.16%  '
    'E18' = '  +3.27%  '
    'E19' = '  +0.15%  '
    'E20' = '  +7.54%  '
    'D21' = '2.047.57'
    'E21' = '  +4.97%  '
    'E22' = '  +3.25%  '
    'E23' = '  +4.16%  '
    'E24' = '  +5.16%  '
    'E25' = '  +3.39%  '
    'E26' = '  +23.62%  '
    'E27' = '  +8.10%  '
    'E28' = '  +6.48%  '
    'E29' = '  +2.85%  '
    'E30' = '  +4.12%  '
    'E31' = '  +5.13%  '
    'E32' = '  +3.07%  '
    'E33' = '  +10.66%  '
    'E34' = '  +8.74%  '
    'E35' = '  +4.60%  '
    'E36' = '  +8.78%  '
    'E37' = '  +12.17%  '
    'E38' = '  +13.03%  '
    'E39' = '  +3.51%  '
    'B40' = 'FraxShare'
    'C40' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'E40' = '  +8.03%  '
    'B41' = 'VeChain'
    'C41' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'E41' = '  +6.63%  '
    'E42' = '  +0.17%  '
    'E43' = '  +6.58%  '
    'E44' = '  -0.69%  '
    'E45' = '  +5.34%  '
    'E46' = '  +5.81%  '
    'E47' = '  +2.70%  '
    'E48' = '  +2.33%  '
    'E49' = '  +4.79%  '
    'B50' = 'NEARProtocol'
    'C50' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E50' = '  +5.83%  '
    'B51' = 'Decentraland'
    'C51' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'E51' = '  +7.95%  '
}

foreach ($addr in $safeValues.Keys) {
    $ws.Range($addr).Value = $safeValues[$addr]
}

# Price values that look like plain numbers (e.g. "0.9999", "248.43") would be
# auto-converted to numeric cells by plain assignment, losing their original
# text formatting. Force them to text by temporarily switching the cell to the
# "@" (Text) number format before assigning, then restore the "Normal" cell
# style so no stray formatting is left behind.
$forceTextValues = [ordered]@{
    'D4' = '0.9999'
    'D5' = '248.43'
    'D6' = '1.000'
    'D7' = '0.4955'
    'D8' = '0.2778'
    'D9' = '0.06416'
    'D11' = '16.75'
    'D12' = '0.07033'
    'D13' = '0.6479'
    'D14' = '84.03'
    'D15' = '4.692'
    'D17' = '0.9998'
    'D18' = '0.000007365'
    'D19' = '0.9999'
    'D20' = '12.26'
    'D22' = '4.562'
    'D23' = '8.918'
    'D24' = '5.333'
    'D25' = '142.40'
    'D26' = '131.81'
    'D27' = '16.50'
    'D28' = '1.884'
    'D29' = '1.409'
    'D30' = '4.155'
    'D31' = '0.08327'
    'D32' = '3.790'
    'D33' = '0.04957'
    'D34' = '1.090'
    'D35' = '2.716'
    'D36' = '0.6720'
    'D37' = '2.261'
    'D38' = '2.760'
    'D39' = '0.9576'
    'D40' = '6.052'
    'D41' = '0.01591'
    'D42' = '0.9999'
    'D43' = '0.4080'
    'D44' = '99.15'
    'D45' = '7.187'
    'D46' = '0.1219'
    'D47' = '0.05519'
    'D48' = '8.086'
    'D49' = '31.49'
    'D50' = '1.306'
    'D51' = '0.3619'
}

foreach ($addr in $forceTextValues.Keys) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $forceTextValues[$addr]
    $ws.Range($addr).Style = "Normal"
}
